$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the last two existing rows (old rows 7 and 8) down to rows 9 and 10
# to make room for the two new rows describing the "general" device type.
$ws.Range("A9").Value2 = "10.9.106.11"
$ws.Range("B9").Value2 = "ex3400"
$ws.Range("C9").Value2 = "oren-flr1sw-B1"
$ws.Range("D9").Value2 = "ex3400"

$ws.Range("A10").Value2 = "192.168.20.243"
$ws.Range("B10").Value2 = "home"
$ws.Range("C10").Value2 = "home-vsrx"
$ws.Range("D10").Value2 = "vsrx"

# Row 6: new general device - pt-router-guest
$ws.Range("A6").Value2 = "10.9.106.46"
$ws.Range("B6").Value2 = "junos"
$ws.Range("C6").Value2 = "pt-router-guest"
$ws.Range("D6").Value2 = "general"

# Row 7: new general device - pt-router-ext1
$ws.Range("A7").Value2 = "10.10.106.1"
$ws.Range("B7").Value2 = "junos"
$ws.Range("C7").Value2 = "pt-router-ext1"
$ws.Range("D7").Value2 = "general"

# Row 8: new task that handles the general device type
$ws.Range("A8").Value2 = "10.10.106.1"
$ws.Range("B8").Value2 = "testing"
$ws.Range("C8").Value2 = "pt-router-ext1"
$ws.Range("D8").Value2 = "general"

[void]$ws.Range("B17").Select()
